$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Pilot Cohort")
$ws1.Columns.Item(14).ColumnWidth = 20.6640625
